# UserData.xlsx: add IsAdmin / IsUser boolean columns (F, G) to the
# UserDataSet sheet, driven by whether the row's login email is the
# admin account (admin@admin.com).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row
$ws.Cells.Item(1, 6).Value = "IsAdmin"
$ws.Cells.Item(1, 7).Value = "IsUser"
$ws.Range("F1:G1").HorizontalAlignment = -4108

# Rows whose Email column (B) is the admin account get IsAdmin = TRUE.
$adminRows = @(12, 14)

for ($r = 2; $r -le 16; $r++) {
    if ($adminRows -contains $r) {
        $ws.Cells.Item($r, 6).Value = $true
    } else {
        $ws.Cells.Item($r, 6).Value = $false
    }
    $ws.Cells.Item($r, 7).Value = $true
}

# IsUser column (G) is centered like the other data columns.
$ws.Range("G2:G16").HorizontalAlignment = -4108

# Leave selection where the author's last interaction landed.
$ws.Range("I13").Select()
